# Natmi following Dr Hou advice
# Update the Artn-Ret LR-pair output sheet with recalculated NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=2.089056333333333; H=6.267169; I=0.7196603919224289; J=0.719660391922429;
           K=3; M=0.8229576666666668; N=2.468873; O=0.03362764644735265; P=0.03362764644735265;
           Q=1.719204925615222; R=15.472844330537; S=0.02420048522173068; T=0.02420048522173069 }
    3 = @{ E=3; G=2.089056333333333; H=6.267169; I=0.7196603919224289; J=0.719660391922429;
           K=3; M=23.10177466666667; N=69.305324; O=0.9439833204831614; P=0.9439833204831615;
           Q=48.26090867863955; R=434.348178107756; S=0.6793474063871477; T=0.679347406387148 }
    4 = @{ E=3; G=2.089056333333333; H=6.267169; I=0.7196603919224289; J=0.719660391922429;
           K=3; M=0.5479189999999999; N=1.643757; O=0.02238903306948597; P=0.02238903306948598;
           Q=1.144633657103666; R=10.301702913933; S=0.01611250031355049; T=0.0161125003135505 }
    5 = @{ E=3; G=0.8137799999999999; H=2.44134; I=0.280339608077571; J=0.280339608077571;
           K=3; M=0.8229576666666668; N=2.468873; O=0.03362764644735265; P=0.03362764644735265;
           Q=0.66970648998; R=6.027358409820001; S=0.009427161225621965; T=0.009427161225621965 }
    6 = @{ E=3; G=0.8137799999999999; H=2.44134; I=0.280339608077571; J=0.280339608077571;
           K=3; M=23.10177466666667; N=69.305324; O=0.9439833204831614; P=0.9439833204831615;
           Q=18.79976218824; R=169.19785969416; S=0.2646359140960136; T=0.2646359140960136 }
    7 = @{ E=3; G=0.8137799999999999; H=2.44134; I=0.280339608077571; J=0.280339608077571;
           K=3; M=0.5479189999999999; N=1.643757; O=0.02238903306948597; P=0.02238903306948598;
           Q=0.4458855238199999; R=4.01296971438; S=0.006276532755935474; T=0.006276532755935475 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
